# The document has a "first page" header/footer pair and a "default"
# (all other pages) header/footer pair, each of which carries an inline
# logo picture:
#   - BTec logo   (JPEG, physically stored as word/media/image2.jpg) sits
#     in both headers.
#   - Pearson logo (PNG, physically stored as word/media/image1.png) sits
#     in both footers.
#
# Their <wp:docPr name="..."> (and the mirrored <pic:cNvPr name="...">)
# labels were mis-keyed to the *other* logo's expected filename; this
# edit renames them to the correct-looking counterpart filenames:
#   BTec header logos:   image2.jpg -> image1.jpg
#   Pearson footer logos: image1.png -> image2.png

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1 (the "default" header/footer used on all
# pages other than a distinct first page)
$wdHeaderFooterPrimary = 1
# wdHeaderFooterFirstPage = 2 (the "first page" header/footer)
$wdHeaderFooterFirstPage = 2

# --- Headers (BTec_Logo-Orange, image2.jpg -> image1.jpg) ---
$headerDefault = $sec.Headers.Item($wdHeaderFooterPrimary)
$headerDefault.Range.InlineShapes.Item(1).Name = "image1.jpg"

$headerFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
$headerFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"

# --- Footers (Pearson logo, image1.png -> image2.png) ---
$footerDefault = $sec.Footers.Item($wdHeaderFooterPrimary)
$footerDefault.Range.InlineShapes.Item(1).Name = "image2.png"

$footerFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
$footerFirst.Range.InlineShapes.Item(1).Name = "image2.png"
